# Applies the FFXIV Leve-profit market-data refresh to all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each data row below is:
#   SheetName|RowNumber|ColumnLetter|NewValue
# NewValue of "CLEAR" means the cell's content is removed entirely
# (matches source rows where the HQ/NQ price lookup returned no data).

$wb = $excel.ActiveWorkbook

$colMap = @{
    "G" = 7
    "H" = 8
    "I" = 9
    "J" = 10
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
}

$changes = @"
ALC|6|H|262.5
ALC|6|I|222.81818
ALC|6|K|668.4545400000001
ALC|6|M|-556.4545400000001
ALC|17|H|3015.8333
ALC|17|J|3119
ALC|17|L|9357
ALC|17|N|-9693
ALC|19|H|1001.9091
ALC|19|J|686.2857
ALC|19|L|686.2857
ALC|19|N|-1036.2857
ALC|32|H|3168.4
ALC|32|I|2621.5
ALC|32|J|3533
ALC|32|K|2621.5
ALC|32|L|3533
ALC|32|M|-2295.5
ALC|32|N|-4185
ALC|40|H|3624.25
ALC|40|I|3248.5
ALC|40|K|3248.5
ALC|40|M|-3073.5
ALC|80|H|51193.082
ALC|80|I|200667.33
ALC|80|J|1368.3334
ALC|80|K|602001.99
ALC|80|L|4105.0002
ALC|80|M|-601003.99
ALC|80|N|-6101.0002
ALC|83|H|51193.082
ALC|83|I|200667.33
ALC|83|J|1368.3334
ALC|83|K|1806005.97
ALC|83|L|12315.0006
ALC|83|M|-1801013.97
ALC|83|N|-22299.0006
ALC|116|H|5791.5
ALC|116|I|5700.1665
ALC|116|K|5700.1665
ALC|116|M|-2258.1665
ARM|61|M|CLEAR
ARM|61|H|3574.5
ARM|61|I|0
ARM|61|K|0
ARM|88|H|1933.6666
ARM|88|J|2586.3333
ARM|88|L|2586.3333
ARM|88|N|-3398.3333
ARM|91|H|1933.6666
ARM|91|J|2586.3333
ARM|91|L|2586.3333
ARM|91|N|-5394.3333
ARM|97|H|9944.923000000001
ARM|97|I|11898.223
ARM|97|J|5550
ARM|97|K|11898.223
ARM|97|L|5550
ARM|97|M|-11402.223
ARM|97|N|-6542
ARM|136|M|CLEAR
ARM|136|H|3574.5
ARM|136|I|0
ARM|136|K|0
BSM|86|M|CLEAR
BSM|86|H|2872.25
BSM|86|I|0
BSM|86|J|2872.25
BSM|86|K|0
BSM|86|L|2872.25
BSM|86|N|-5118.25
BSM|89|M|CLEAR
BSM|89|H|2872.25
BSM|89|I|0
BSM|89|J|2872.25
BSM|89|K|0
BSM|89|L|14361.25
BSM|89|N|-25593.25
BSM|95|H|48999
BSM|95|J|48999
BSM|95|L|48999
BSM|95|N|-54491
BSM|134|H|1798.4166
BSM|134|I|1707.6818
BSM|134|K|5123.0454
BSM|134|M|-2588.0454
CRP|7|H|76.10526
CRP|7|I|67.46154
CRP|7|K|67.46154
CRP|7|M|45.53846
CRP|22|H|1112.25
CRP|22|I|306.57144
CRP|22|J|1738.8889
CRP|22|K|306.57144
CRP|22|L|1738.8889
CRP|22|M|43.42856
CRP|22|N|-2438.8889
CRP|86|H|19363.814
CRP|86|I|29102.562
CRP|86|J|5198.364
CRP|86|K|29102.562
CRP|86|L|5198.364
CRP|86|M|-27979.562
CRP|86|N|-7444.364
CRP|89|H|19363.814
CRP|89|I|29102.562
CRP|89|J|5198.364
CRP|89|K|145512.81
CRP|89|L|25991.82
CRP|89|M|-139896.81
CRP|89|N|-37223.82
CRP|132|H|3101.2
CRP|132|I|2836.7646
CRP|132|J|4599.6665
CRP|132|K|8510.293799999999
CRP|132|L|13798.9995
CRP|132|M|-5980.293799999999
CRP|132|N|-18858.9995
CRP|134|H|167826
CRP|134|I|201195.6
CRP|134|K|603586.8
CRP|134|M|-601051.8
CRP|141|H|387460.38
CRP|141|J|433383.28
CRP|141|L|433383.28
CRP|141|N|-443743.28
CUL|4|H|1245901.2
CUL|4|I|869858.9399999999
CUL|4|J|3000765.8
CUL|4|K|2609576.82
CUL|4|L|9002297.399999999
CUL|4|M|-2609464.82
CUL|4|N|-9002521.399999999
CUL|122|H|447.82352
CUL|122|I|358.1
CUL|122|J|576
CUL|122|K|3222.9
CUL|122|L|5184
CUL|122|M|-772.9000000000001
CUL|122|N|-10084
GSM|21|H|14999
GSM|21|J|14999
GSM|21|L|14999
GSM|21|N|-15345
GSM|30|H|14999
GSM|30|J|14999
GSM|30|L|14999
GSM|30|N|-15209
GSM|39|H|49986.5
GSM|39|J|49986.5
GSM|39|L|49986.5
GSM|39|N|-51050.5
GSM|70|H|6599
GSM|70|I|4799
GSM|70|K|4799
GSM|70|M|-4529
GSM|73|H|6599
GSM|73|I|4799
GSM|73|K|4799
GSM|73|M|-3863
GSM|102|H|3516.5454
GSM|102|I|2520.3333
GSM|102|K|2520.3333
GSM|102|M|-898.3332999999998
GSM|132|H|87457.086
GSM|132|I|101948.6
GSM|132|J|14999.5
GSM|132|K|305845.8
GSM|132|L|44998.5
GSM|132|M|-303315.8
GSM|132|N|-50058.5
GSM|136|H|40122.453
GSM|136|J|40122.453
GSM|136|L|120367.359
GSM|136|N|-125467.359
LTW|16|H|831.875
LTW|16|I|736.4286
LTW|16|J|1500
LTW|16|K|736.4286
LTW|16|L|1500
LTW|16|M|-566.4286
LTW|16|N|-1840
LTW|22|H|33332.324
LTW|22|J|2432.7778
LTW|22|L|2432.7778
LTW|22|N|-3022.7778
LTW|27|H|33332.324
LTW|27|J|2432.7778
LTW|27|L|2432.7778
LTW|27|N|-2646.7778
LTW|46|H|12372.363
LTW|46|I|31066.334
LTW|46|K|31066.334
LTW|46|M|-30878.334
LTW|55|H|1275.091
LTW|55|I|1211.9333
LTW|55|K|1211.9333
LTW|55|M|-1038.9333
LTW|93|H|1628.1666
LTW|93|I|817
LTW|93|K|817
LTW|93|M|431
WVR|62|H|134242.22
WVR|62|I|4733.6665
WVR|62|K|4733.6665
WVR|62|M|-4109.6665
WVR|65|H|134242.22
WVR|65|I|4733.6665
WVR|65|K|23668.3325
WVR|65|M|-20548.3325
WVR|126|H|514999.34
WVR|126|I|514999.34
WVR|126|K|1544998.02
WVR|126|M|-1542528.02
WVR|132|H|76696.75
WVR|132|I|80101.92
WVR|132|J|34699.668
WVR|132|K|240305.76
WVR|132|L|104099.004
WVR|132|M|-237775.76
WVR|132|N|-109159.004
WVR|136|H|6080.8887
WVR|136|I|6885.5
WVR|136|J|3782
WVR|136|K|20656.5
WVR|136|L|11346
WVR|136|M|-18106.5
WVR|136|N|-16446
"@

$sheetCache = @{}
$appliedCount = 0

$lines = $changes -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }

    $parts = $line -split '\|'
    $sheetName = $parts[0]
    $rowNum = [int]$parts[1]
    $colLetter = $parts[2]
    $valToken = $parts[3]
    $colNum = $colMap[$colLetter]

    if (-not $sheetCache.ContainsKey($sheetName)) {
        $sheetCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $sheetCache[$sheetName]

    $cell = $ws.Cells.Item($rowNum, $colNum)
    if ($valToken -eq "CLEAR") {
        $cell.ClearContents()
    } else {
        $cell.Value = [double]$valToken
    }
    $appliedCount = $appliedCount + 1
}

Write-Output "Applied $appliedCount cell updates across $($sheetCache.Count) sheets."
